# Update "paises.xlsx" (Pais sheet) with the latest COVID country/provincia
# numbers and refresh the "datos actualizados" timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Refresh "last updated" timestamp (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 19 de Mayo de 2020 a las 20:05"

# --- Update per-country numeric figures (columns B:H) ---
# Each entry: row -> Casos totales, Nuevos casos, Casos activos, Recuperados,
#             Casos criticos, Muertes hoy, Muertes
$updates = @{
    4   = @(1558498, 8204, 360058, 1105920, 0, 539, 92520)
    11  = @(177728, 439, 155700, 13860, 0, 45, 8168)
    14  = @(106446, 6118, 40865, 62280, 0, 145, 3301)
    17  = @(79070, 998, 39951, 33210, 0, 67, 5909)
    44  = @(13484, 720, 3742, 9083, 0, 14, 659)
    45  = @(13223, 498, 6613, 6169, 0, 7, 441)
    46  = @(12942, 224, 2843, 9262, 0, 6, 837)
    68  = @(3958, 11, 3718, 131, 0, 2, 109)
    105 = @(1023, 31, 569, 445, 0, 0, 9)
    139 = @(335, 7, 85, 247, 0, 0, 3)
    140 = @(335, 0, 300, 11, 0, 0, 24)
    141 = @(332, 0, 322, 0, 0, 0, 10)
    142 = @(330, 0, 106, 212, 0, 0, 12)
}

$cols = @("B", "C", "D", "E", "F", "G", "H")

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$row").Value = $vals[$i]
    }
}
